$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "Main"

# Update A1 value
$ws.Range("A1").Value = 9

# Fill column A (item codes) first
$ws.Range("A2").Value = "a1"
$ws.Range("A3").Value = "a2"
$ws.Range("A4").Value = "a3"
$ws.Range("A5").Value = "b1"
$ws.Range("A6").Value = "b2"
$ws.Range("A7").Value = "b3"
$ws.Range("A8").Value = "c1"
$ws.Range("A9").Value = "c2"
$ws.Range("A10").Value = "c3"

# Fill column B (status) second
$ws.Range("B2").Value = "o"
$ws.Range("B3").Value = "o"
$ws.Range("B4").Value = "o"
$ws.Range("B5").Value = "e"
$ws.Range("B6").Value = "o"
$ws.Range("B7").Value = "o"
$ws.Range("B8").Value = "e"
$ws.Range("B9").Value = "e"
$ws.Range("B10").Value = "e"

# Update the active selection shown in the sheet view
$ws.Range("B11").Select()
